$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# wdFindContinue = 1, wdReplaceAll = 2
$find.Execute(
    "môžete pozorovať súhvezdie ozvezdje Pegaz 2022: 8.-17. oktober, 7.-16. november,",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "2022: Datumi kampanje za opazovanje ozvezdje Pegaz: 8.-17. oktober, 7.-16. november,",
    2
)
